# "Corrected spoofed and non-spoofed folder and updated contents of folders"
#
# The sheet physically behind the "Non-spoofed" tab is actually the spoofed
# results and vice-versa, so the two tab names need to be swapped (the
# underlying sheets/data stay where they are - only the labels move).
# The (relabeled) "Spoofed" tab then becomes the active/selected tab, scrolled
# down to the bottom table, and a recomputed statistic is written back.

$wb = $excel.ActiveWorkbook

$sNonSpoofed = $wb.Worksheets.Item("Non-spoofed")
$sSpoofed    = $wb.Worksheets.Item("Spoofed")

# Swap the two tab names. Excel won't allow two sheets to share a name even
# momentarily, so stage the rename through a temporary, certainly-unique name.
$tempName = "__tmp_rename__"
$sNonSpoofed.Name = $tempName
$sSpoofed.Name    = "Non-spoofed"
$sNonSpoofed.Name = "Spoofed"

# $sNonSpoofed now carries the "Spoofed" label (rId2/sheetId 2) - make it the
# active sheet/tab and move the on-sheet selection to where the author left it.
[void]$sNonSpoofed.Activate()
[void]$sNonSpoofed.Range("J57").Select()

# Scroll the view so row 42 is pinned at the top of the window, matching the
# author's saved scroll position.
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1

# The underlying statistic in T68 was recomputed (a last-digit float nudge);
# O15 references T68 directly (=T68) so it recalculates to the same new value
# automatically, without needing to touch its formula.
$sNonSpoofed.Range("T68").Value = 0.93367811667623302
